# Auto-generated Excel COM-interop script applying the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.763.75"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "'2.150.02"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'352.29"
$ws.Range("E5").Value = "  +5.41%  "
$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "'0.5288"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").Value = "'0.4573"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("D9").Value = "'54.27"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +3.38%  "
$ws.Range("D11").Value = "'1.187"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E12").Value = "  +3.76%  "
$ws.Range("D13").Value = "'2.150.88"
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "'6.923"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "'102.54"
$ws.Range("E16").Value = "  +6.11%  "
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'0.06731"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "'19.65"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "'6.378"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").Value = "'30.854.44"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "'12.90"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").Value = "'2.399"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").Value = "'2.374.21"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'22.70"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").Value = "'2.639"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("D29").Value = "'165.32"
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("D30").Value = "'137.05"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("D31").Value = "'1.226"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("E33").Value = "  +3.31%  "
$ws.Range("D34").Value = "'6.418"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "'4.019"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").Value = "'6.164"
$ws.Range("E36").Value = "  +6.15%  "
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").Value = "'0.02663"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").Value = "'0.06939"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D40").Value = "'0.2345"
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").Value = "'0.7004"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("D43").Value = "'1.277"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").Value = "'14.76"
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").Value = "'0.6500"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000370"
$ws.Range("E47").Value = "  +5.42%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.756"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").Value = "'1.263"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").Value = "'83.60"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("D51").Value = "'0.07328"
$ws.Range("E51").Value = "  +2.56%  "
